$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(6,7,9,10,17,20,25,26,29,31,34,39,44,47,52,57,58,62,69,70,76,81,86,89,92,93,94,98,101,102,106,111,114,116,118,122,124,127,128,132,135,138,139,141,142,143,148,151,153,155,157,159,160,161,162,169,170,176,180,182,183,190,191,192)

foreach ($r in $rows) {
    $ws.Range("G$r`:N$r").Value = "N/A"
}
